$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.087.60'
$ws.Range('E2').Value = '  +0.80%  '
$ws.Range('D3').Value = '2.106.40'
$ws.Range('E3').Value = '  +10.26%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = "'253.29"
$ws.Range('E5').Value = '  +1.62%  '
$ws.Range('D6').Value = "'0.663"
$ws.Range('E6').Value = '  -4.88%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = "'50.25"
$ws.Range('E8').Value = '  +7.33%  '
$ws.Range('D9').Value = "'60.93"
$ws.Range('E9').Value = '  +5.10%  '
$ws.Range('E10').Value = '  +0.95%  '
$ws.Range('D11').Value = "'0.0748"
$ws.Range('E11').Value = '  -1.34%  '
$ws.Range('E12').Value = '  +6.92%  '
$ws.Range('D13').Value = "'14.81"
$ws.Range('E13').Value = '  +0.50%  '
$ws.Range('D14').Value = '2.406.25'
$ws.Range('E14').Value = '  +10.01%  '
$ws.Range('D15').Value = "'0.841"
$ws.Range('E15').Value = '  +3.41%  '
$ws.Range('D16').Value = '2.112.47'
$ws.Range('E16').Value = '  +10.10%  '
$ws.Range('E17').Value = '  +1.82%  '
$ws.Range('D18').Value = '36.969.31'
$ws.Range('E18').Value = '  +0.63%  '
$ws.Range('D19').Value = "'73.42"
$ws.Range('E19').Value = '  -1.28%  '
$ws.Range('D20').Value = '0.0₃0829'
$ws.Range('E20').Value = '  -2.89%  '
$ws.Range('D21').Value = "'13.37"
$ws.Range('E21').Value = '  -2.07%  '
$ws.Range('D22').Value = "'241.80"
$ws.Range('E22').Value = '  -3.57%  '
$ws.Range('E23').Value = '  +3.61%  '
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('D25').Value = "'2.53"
$ws.Range('E25').Value = '  +1.49%  '
$ws.Range('D26').Value = "'170.60"
$ws.Range('E26').Value = '  +2.11%  '
$ws.Range('D27').Value = "'9.40"
$ws.Range('E27').Value = '  +7.20%  '
$ws.Range('D28').Value = "'21.25"
$ws.Range('E28').Value = '  +13.79%  '
$ws.Range('D29').Value = "'2.02"
$ws.Range('E29').Value = '  -7.30%  '
$ws.Range('D30').Value = "'26.55"
$ws.Range('E30').Value = '  +38.27%  '
$ws.Range('E31').Value = '  -3.98%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = "'4.54"
$ws.Range('E32').Value = '  -1.40%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = "'1.05"
$ws.Range('E33').Value = '  +21.28%  '
$ws.Range('E34').Value = '  +1.13%  '
$ws.Range('D35').Value = "'0.0928"
$ws.Range('E35').Value = '  +5.35%  '
$ws.Range('D36').Value = "'2.42"
$ws.Range('E36').Value = '  +22.31%  '
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('D38').Value = "'4.13"
$ws.Range('E38').Value = '  -4.16%  '
$ws.Range('E39').Value = '  -6.50%  '
$ws.Range('E40').Value = '  -9.28%  '
$ws.Range('E41').Value = '  -0.49%  '
$ws.Range('E42').Value = '  +7.50%  '
$ws.Range('D43').Value = "'98.89"
$ws.Range('E43').Value = '  -5.38%  '
$ws.Range('D44').Value = "'17.04"
$ws.Range('E44').Value = '  -4.45%  '
$ws.Range('D45').Value = "'2.79"
$ws.Range('E45').Value = '  -2.03%  '
$ws.Range('D46').Value = '1.355.88'
$ws.Range('E46').Value = '  +0.23%  '
$ws.Range('D47').Value = "'0.0859"
$ws.Range('E47').Value = '  +4.98%  '
$ws.Range('B48').Value = 'MXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D48').Value = "'2.95"
$ws.Range('E48').Value = '  +4.75%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Value = "'7.11"
$ws.Range('E49').Value = '  +11.26%  '
$ws.Range('D50').Value = '2.290.72'
$ws.Range('E50').Value = '  +9.73%  '
$ws.Range('D51').Value = "'2.28"
$ws.Range('E51').Value = '  -3.96%  '
